$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: insert a new value (0.16, unstyled) into B7, and shift the
#     existing B7 value/style (2.5E-2, style s=1) into the new C7 cell. ---
$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("B7").Style = "Normal"
$ws.Range("B7").Value = 0.16

# --- Row 8: give it an explicit custom row height. ---
$ws.Rows.Item(8).RowHeight = 15.6

# --- Rows 19-22: insert a new value (1, unstyled) into column B, shifting
#     the existing B value/style into the new C cell. ---
$rows = 19, 20, 21, 22
foreach ($r in $rows) {
    $bCell = "B" + $r
    $cCell = "C" + $r
    $ws.Range($bCell).Copy()
    $ws.Range($cCell).PasteSpecial()
    $ws.Range($cCell).PasteSpecial(-4122)
    $ws.Range($bCell).Style = "Normal"
    $ws.Range($bCell).Value = 1
}

# --- Update the saved selection / active cell. ---
$ws.Range("B23").Select() | Out-Null
